$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows (prices/volumes/dates) were reshuffled between
# rows: row2<->row5, row4<->row8, row6<->row7 (row3 stays the same).
# Apply the new values for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).

$rowValues = @{
    2 = @{ D = 45092; J = 210; K = 10000; L = 11000; M = 10714; P = 595 }
    4 = @{ D = 45215; J = 200; K = 11000; L = 12000; M = 11500; P = 639 }
    5 = @{ D = 45175; J = 250; K = 11000; L = 12000; M = 11500; P = 639 }
    6 = @{ D = 44714; J = 80;  K = 9000;  L = 10000; M = 9500;  P = 528 }
    7 = @{ D = 45205; J = 200; K = 11000; L = 12000; M = 11500; P = 639 }
    8 = @{ D = 44792; J = 160; K = 9000;  L = 10000; M = 9500;  P = 528 }
}

foreach ($row in $rowValues.Keys) {
    $vals = $rowValues[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
